# Add a copy of the "Release EDI" data as a second worksheet named
# "Release_EDI" (underscore), with a matching defined name, per the
# "updates catch, release, recapture tables from grant 12-15-2022" commit.

$wb = $excel.ActiveWorkbook

# Grid of values for the new sheet: header row + 20 data rows, columns A:P.
# $null marks a cell that is intentionally left blank (no cell written at
# all, matching the source data which omits NA values entirely).
$rows = @(
  ,@("projectDescriptionID", "releaseID", "commonName", "markedRun", "markedLifeStage", "markedFishOrigin", "sourceOfFishSite", "releaseSite", "releaseSubSite", "nReleased", "releaseTime", "testDays", "appliedMarkType", "appliedMarkColor", "appliedMarkPosition", "appliedMarkCode")
  ,@(0, 0, "Not applicable (n/a)", "Not recorded", "Not recorded", "Unknown", "Not applicable", "Not applicable", $null, 0, 1, 0, $null, $null, $null, $null)
  ,@(0, 255, "Not applicable (n/a)", "Not recorded", "Not recorded", "Unknown", "Not applicable", "Not applicable", $null, 0, 1, 0, $null, $null, $null, $null)
  ,@(11, 0, "Not applicable (n/a)", "Not applicable (n/a)", "Not applicable (n/a)", "Not applicable (n/a)", "Not applicable", "Not applicable", $null, $null, $null, 0, $null, $null, $null, $null)
  ,@(11, 255, "Not applicable (n/a)", "Not applicable (n/a)", "Not applicable (n/a)", "Not applicable (n/a)", "Not applicable", "Not applicable", $null, $null, $null, 0, $null, $null, $null, $null)
  ,@(11, 256, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, $null, 100, 44202.5750578704, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 257, "Chinook salmon", "Spring", $null, "Natural", $null, $null, $null, 100, 44202.5208680556, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 258, "Chinook salmon", "Spring", $null, "Natural", $null, $null, $null, 100, 44481.5832060185, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 259, "Not applicable (n/a)", $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  ,@(11, 260, "Not applicable (n/a)", $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  ,@(11, 261, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, $null, 100, 44202.5208680556, 6, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 262, "Not applicable (n/a)", $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
  ,@(11, 263, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, "n/a", 249, 44210.5213657407, 6, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 264, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, "n/a", 109, 44237.5002893519, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 265, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, "n/a", 349, 44251.5209143519, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 266, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, "n/a", 197, 44265.5004166667, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 267, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, "n/a", 160, 44294.5418171296, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 268, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, "n/a", 500, 44301.5418865741, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 269, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, "n/a", 250, 44600.5106365741, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 270, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, "n/a", 249, 44607.5523032407, 7, "Pigment / dye", "Brown", "Whole body", $null)
  ,@(11, 271, "Chinook salmon", "Spring", "Juvenile", "Natural", $null, $null, "n/a", 185, 44628.555775463, 7, "Pigment / dye", "Brown", "Whole body", $null)
)

# New sheet goes right after the existing "Release EDI" tab.
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $firstSheet)
$ws.Name = "Release_EDI"

$dateCol = 11  # column K = releaseTime

for ($r = 0; $r -lt $rows.Count; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $val = $row[$c]
        if ($null -ne $val) {
            $cell = $ws.Cells.Item($r + 1, $c + 1)
            $cell.Value = $val
            if ($c -eq ($dateCol - 1) -and $r -gt 0) {
                $cell.NumberFormat = "m/d/yyyy"
            }
        }
    }
}

# Match the source sheet's view: active cell C5 selected, first tab shown.
$ws.Range("C5").Select()

# Defined name spanning the full written range of the new sheet.
$wb.Names.Add("Release_EDI", "='Release_EDI'!`$A`$1:`$P`$21")
